# leetcode problem 933; queue problem using deque() class
# Adds two new rows (6 & 7) to the "July" sheet:
#   row 6 - "maximum average subarray" (easy, #643, sliding window)
#   row 7 - "Maximum Number of Vowels in a Substring of Given Length"
#           (medium, #1456, sliding window, with submission link)
# Also reformats the sheet: centers the data columns, wraps / widens the
# notes column, and bumps the row heights to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July")

# ---------------------------------------------------------------------
# 1. New cell values, written in the same order the original author
#    would have typed them (this also controls shared-string ordering).
# ---------------------------------------------------------------------

# Row 6 - maximum average subarray
$ws.Range("A6").Value = "easy"
$ws.Range("B6").Value = 643
$ws.Range("C6").Value = "maximum average subarray "

# Row 7 - Maximum Number of Vowels in a Substring of Given Length
$ws.Range("A7").Value = "medium"
$ws.Range("B7").Value = 1456
$ws.Range("D7").Value = "sliding window"
$ws.Range("C7").Value = "Maximum Number of Vowels in a Substring of Given Length"

$ws.Hyperlinks.Add($ws.Range("E7"), "https://leetcode.com/problems/maximum-number-of-vowels-in-a-substring-of-given-length/submissions/1321986754?envType=study-plan-v2&envId=leetcode-75")

# ---------------------------------------------------------------------
# 2. Formatting
# ---------------------------------------------------------------------

# Widen the notes column and let it wrap.
$ws.Columns.Item(3).ColumnWidth = 32
$ws.Range("C1:C7").WrapText = $true

# Center the Tag / Problem number / notes / Link columns for all data rows.
$ws.Range("A2:B7").HorizontalAlignment = -4108
$ws.Range("D2:D7").HorizontalAlignment = -4108
$ws.Range("C2:C7").HorizontalAlignment = -4108

# Row heights: the header + existing rows grow slightly, the new
# two-line submission-link row is taller.
$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 17
$ws.Rows.Item(3).RowHeight = 17
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 34

# ---------------------------------------------------------------------
# 3. Selection - the author left the cursor on the first empty row.
# ---------------------------------------------------------------------
$ws.Range("A8").Select()
